# ---------------------------------------------------------------------
# cs-en-us-046pct.xlsx weekly refresh: new crime data collected
# ---------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text: "Volume 32 Number 4" -> "Volume 32 Number 5" ---------
$hdr = $ws.Range("A8")
$hdr.Characters(21, 1).Text = "5"

# --- Header text: reporting week dates ----------------------------------
$week = $ws.Range("C9")
$week.Characters(27, 9).Text = "1/27/2025"
$week.Characters(47, 9).Text = "2/2/2025"

# --- Cells that keep their existing number format, only the value moves 
$sameStyleValues = @{
    "J14" = 2
    "F15" = 3
    "G15" = 1
    "H15" = 200
    "I15" = 3
    "K15" = 50
    "L15" = 50
    "C16" = 4
    "D16" = 6
    "E16" = -33.333333333333
    "F16" = 26
    "G16" = 26
    "H16" = 0
    "I16" = 28
    "J16" = 32
    "K16" = -12.5
    "L16" = -9.677419354838
    "M16" = -17.647058823529
    "N16" = -87.214611872146
    "D17" = 10
    "E17" = 20
    "G17" = 50
    "H17" = 6
    "I17" = 59
    "J17" = 64
    "K17" = -7.8125
    "L17" = 7.272727272727
    "M17" = 103.448275862069
    "N17" = -37.234042553191
    "C18" = 1
    "D18" = 3
    "E18" = -66.666666666666
    "F18" = 5
    "H18" = -66.666666666666
    "I18" = 7
    "J18" = 18
    "K18" = -61.111111111111
    "L18" = -53.333333333333
    "M18" = -73.076923076923
    "N18" = -96.682464454976
    "C19" = 13
    "D19" = 15
    "E19" = -13.333333333333
    "F19" = 47
    "G19" = 44
    "H19" = 6.818181818181
    "I19" = 51
    "J19" = 54
    "K19" = -5.555555555555
    "L19" = -1.923076923076
    "M19" = 142.857142857143
    "N19" = -40
    "C20" = 4
    "D20" = 2
    "E20" = 100
    "F20" = 21
    "G20" = 14
    "H20" = 50
    "I20" = 22
    "J20" = 15
    "K20" = 46.666666666666
    "L20" = -31.25
    "M20" = 144.444444444444
    "N20" = -73.809523809523
    "C21" = 36
    "D21" = 37
    "E21" = -2.702702702702
    "F21" = 155
    "G21" = 151
    "H21" = 2.649006622516
    "I21" = 170
    "J21" = 187
    "K21" = -9.090909090909
    "L21" = -9.090909090909
    "M21" = 41.666666666666
    "N21" = -75.783475783475
    "G22" = 2
    "J22" = 3
    "F23" = 2
    "I23" = 2
    "L23" = -60
    "M23" = -33.333333333333
    "C24" = 15
    "D24" = 20
    "E24" = -25
    "F24" = 70
    "G24" = 71
    "H24" = -1.408450704225
    "I24" = 83
    "J24" = 78
    "K24" = 6.410256410256
    "L24" = 1.219512195121
    "M24" = -12.631578947368
    "C25" = 2
    "D25" = 3
    "E25" = -33.333333333333
    "F25" = 14
    "G25" = 16
    "H25" = -12.5
    "I25" = 17
    "J25" = 17
    "K25" = 0
    "L25" = -34.615384615384
    "C26" = 10
    "D26" = 32
    "E26" = -68.75
    "F26" = 60
    "G26" = 100
    "H26" = -40
    "I26" = 72
    "J26" = 111
    "K26" = -35.135135135135
    "L26" = -31.428571428571
    "M26" = -22.580645161290
    "F27" = 3
    "G27" = 3
    "H27" = 0
    "I27" = 3
    "K27" = -40
    "L27" = -62.5
    "D28" = 3
    "E28" = -66.666666666666
    "F28" = 5
    "G28" = 9
    "H28" = -44.444444444444
    "I28" = 6
    "J28" = 10
    "K28" = -40
    "L28" = -25
    "N29" = -93.333333333333
    "N30" = -92.307692307692
}
foreach ($addr in $sameStyleValues.Keys) {
    $ws.Range($addr).Value = $sameStyleValues[$addr]
}

# --- Cells that previously held the text placeholder ("0" / "***.*") and
#     now hold a real number -> give them the matching numeric format ----
$numFormat = "#,##0"
$pctFormat = '#,##0.0;"-"#,##0.0'
$toNumeric = @{
    "D14" = @(1, $numFormat)
    "E14" = @(-100, $pctFormat)
    "C15" = @(2, $numFormat)
    "N15" = @(0, $pctFormat)
    "C23" = @(1, $numFormat)
    "D23" = @(1, $numFormat)
    "E23" = @(0, $pctFormat)
    "G23" = @(1, $numFormat)
    "H23" = @(100, $pctFormat)
    "J23" = @(1, $numFormat)
    "K23" = @(100, $pctFormat)
    "C27" = @(2, $numFormat)
}
foreach ($addr in $toNumeric.Keys) {
    $pair = $toNumeric[$addr]
    $ws.Range($addr).Value = $pair[0]
    $ws.Range($addr).NumberFormat = $pair[1]
}

# --- Cells that go the other way: numeric -> back to the text "0" -------
#     placeholder. Force text storage via a Text (@) format, then restore
#     the original display format (borrowed from an untouched sibling),
#     matching how Excel stores a typed "0" in a Text-formatted cell.
$textTemplate = $ws.Range("C14")
$textTemplate.Copy() | Out-Null
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "0"
$textTemplate.Copy() | Out-Null
$ws.Range("C29").PasteSpecial(-4122)
$textTemplate.Copy() | Out-Null
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "0"
$textTemplate.Copy() | Out-Null
$ws.Range("C30").PasteSpecial(-4122)

Write-Host "Applied weekly crime-data refresh."
